$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductionPlan")

# Row 2 (MAT_A / LINE_A)
$ws.Range("G2").Value = 860
$ws.Range("H2").Value = 860
$ws.Range("J2").Value = 817

# Row 3 (MAT_B / LINE_B)
$ws.Range("G3").Value = 112
$ws.Range("H3").Value = 112
$ws.Range("J3").Value = 99
